# "Bausteinsicht hinzugefuegt, einige Ignores und Literaturverzeichnis"
#
# Update the task-status column (E) of the Projektplan sheet to reflect
# progress: several tasks moved to "Erledigt" / "in Bearbeitung", and a
# batch of later tasks are now flagged "ueberfaellig" (a brand-new status
# string that gets added to the shared-string table).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- status column updates -------------------------------------------------
$ws.Range("E5").Value  = "Erledigt"
$ws.Range("E6").Value  = "Erledigt"
$ws.Range("E7").Value  = "Erledigt"
$ws.Range("E8").Value  = "Erledigt"

$ws.Range("E9").Value  = "in Bearbeitung"
$ws.Range("E10").Value = "in Bearbeitung"

$ws.Range("E11").Value = "überfällig"
$ws.Range("E12").Value = "überfällig"
$ws.Range("E13").Value = "überfällig"
$ws.Range("E14").Value = "überfällig"
$ws.Range("E15").Value = "überfällig"
$ws.Range("E16").Value = "überfällig"
$ws.Range("E17").Value = "überfällig"

# --- current selection, as left by the editing session ----------------------
$ws.Range("F17:G17").Select()

# --- page setup: printed at 52% scale, "fit to page" flag left set ---------
$ps = $ws.PageSetup
$ps.Zoom = 52
$ps.FitToPagesWide = $false
$ps.FitToPagesTall = $false

# sheet-level flag toggled alongside the conditional-formatting cleanup
$ws.EnableFormatConditionsCalculation = $false

# --- window position, as saved with the workbook ----------------------------
$win = $excel.Windows.Item(1)
$win.Left = 1780
$win.Top = 760
